$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 updates (Bagmati) ---
$ws.Range("C4").Value = 23
$ws.Range("D4").Value = 234.98999999999998
$ws.Range("E4").Value = 230.06000000000003
$ws.Range("F4").Value = 97.902038384612126
$ws.Range("G4").Value = 77.323295129902874

# --- Row 6 updates (Karnali) ---
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 21.94
$ws.Range("E6").Value = 22.34
$ws.Range("F6").Value = 101.82315405651777
$ws.Range("G6").Value = 7.5084865391725195

# --- Row 7 updates (Sudurpashchim) ---
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").ClearContents()
$ws.Range("G7").Value = 0

# --- Column J: replace formulas with static values (J2:J9) ---
$ws.Range("J2").Value = 131.37105531107591
$ws.Range("J3").Value = 96.067296208716698
$ws.Range("J4").Value = 244.1956981031831
$ws.Range("J5").Value = 139.47241555329879
$ws.Range("J6").Value = 118.71063494395266
$ws.Range("J7").Value = 102.35364119934859
$ws.Range("J8").Value = 106.46186720689109
$ws.Range("J9").Value = 142.18206141637918

# --- Selection change ---
$ws.Range("J2:J9").Select()

$wb.Save()
